$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = 0.638060484538402
$ws.Range("A14").Value = 0.6503828814202762
$ws.Range("A18").Value = 0.8022805061070414
$ws.Range("D2").Value = 1.067944483982438
$ws.Range("E2").Value = 0.00491090676638329
$ws.Range("D3").Value = 1.041063796643669
$ws.Range("E3").Value = 0.04798539450922949
$ws.Range("D4").Value = 0.8724697628914463
$ws.Range("E4").Value = 0.1070374947921918
$ws.Range("D5").Value = 0.7222702532140244
$ws.Range("E5").Value = 0.2768709729398944
$ws.Range("D6").Value = 0.8087891143948103
$ws.Range("E6").Value = 0.2274939700628229
$ws.Range("D7").Value = 0.7491511768755239
$ws.Range("E7").Value = 0.7442763759308628
$ws.Range("D8").Value = 0.8396520785252934
$ws.Range("E8").Value = 0.9589726175047373
$ws.Range("D9").Value = 0.8773066217241797
$ws.Range("E9").Value = 0.9583160867926905
$ws.Range("D10").Value = 1.046104459252619
$ws.Range("E10").Value = 1.305287256944929
$ws.Range("D11").Value = 1.071577964559529
$ws.Range("E11").Value = 1.369601345896547
$ws.Range("D12").Value = 1.085316886986227
$ws.Range("E12").Value = 1.395265567373854
$ws.Range("D13").Value = 1.063824758802943
$ws.Range("E13").Value = 1.363245809946933
$ws.Range("D14").Value = 1.082410418953807
$ws.Range("E14").Value = 1.385093570686751
$ws.Range("D15").Value = 1.08261330774133
$ws.Range("E15").Value = 1.387122359146569
$ws.Range("D16").Value = 1.096662296799546
$ws.Range("E16").Value = 1.406081803184485
$ws.Range("D17").Value = 1.092223243382248
$ws.Range("E17").Value = 1.403742817108704
$ws.Range("D18").Value = 1.100277460873889
$ws.Range("E18").Value = 1.414880450298408
$ws.Range("D19").Value = 1.09617353675141
$ws.Range("E19").Value = 1.409199123762568
$ws.Range("D20").Value = 1.104545772992451
$ws.Range("E20").Value = 1.42079604921987
$ws.Range("D21").Value = 1.09962260465262
$ws.Range("E21").Value = 1.413820027131574
